$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "zip"
$ws.Columns.Item(1).ColumnWidth = 20.5
